# Append: 2026-02-17 13:12 JST
# Update the "取得日時" (acquired datetime) column for all existing data rows
# on the "ランサーズ" sheet from 2026-02-17 12:58:45 to 2026-02-17 13:12:31.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2026-02-17 12:58:45"
$newValue = "2026-02-17 13:12:31"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
